# Re-order rows 5-23 (skipping 13, 15, 19 which are unaffected) of the
# "Artfynd" sheet according to the source commit's data refresh. The edit
# is a pure permutation of entire data rows (every column's content moves
# together) - no individual field values are altered, only which physical
# row the already-existing record occupies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that participate in the shuffle.
$rows = @(5, 6, 7, 8, 9, 10, 11, 12, 14, 16, 17, 18, 20, 21, 22, 23)

# new-row -> old-row the data should come from.
$mapping = @{
    5  = 7
    6  = 5
    7  = 6
    8  = 16
    9  = 10
    10 = 18
    11 = 17
    12 = 23
    14 = 20
    16 = 22
    17 = 8
    18 = 9
    20 = 14
    21 = 12
    22 = 21
    23 = 11
}

# Columns A..AY (1..51) that hold genuinely numeric (t="n") data - these
# must stay numbers when copied.
$numericCols = @(1, 2, 5, 17, 18, 19)
# Columns that hold boolean (t="b") data - must stay booleans.
$boolCols = @(30, 31, 33)

$lastCol = 51

# 1) Snapshot the current content of every affected row (whole A:AY range)
#    before any writes happen, so the permutation can't clobber a row
#    that hasn't been read yet.
$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, $lastCol)).Value2
}

# 2) Write each row's new content from its mapped source snapshot.
foreach ($r in $rows) {
    $src = $mapping[$r]
    # Arrays coming back from Excel's Value2 are 1-based in both
    # dimensions ([1..1, 1..51]) - mutate in place (a fresh
    # New-Object array would be 0-based and mismatch Excel's writer).
    $data = $snapshot[$src]
    for ($c = 1; $c -le $lastCol; $c++) {
        $val = $data[1, $c]
        if ($null -ne $val -and $numericCols -notcontains $c -and $boolCols -notcontains $c) {
            $data[1, $c] = "'" + $val
        }
    }

    $destRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, $lastCol))
    $destRange.Value2 = $data
    # Clear the "quote prefix" style the apostrophe entry implicitly
    # applies, so the cell style stays identical to the untouched cells.
    $destRange.Style = "Normal"
}
